$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused trailing rows (old sheet had 15 data rows, new one has 12)
$ws.Rows("13:15").Delete()

# Overwrite the remaining rows 1-12 with the new data set
$data = @(
    @(5272964724, "4X40Y95214", 32,  "ZC06"),
    @(5272964724, "21MDS2HH07", 9,   "ZC06"),
    @(5272964724, "40AY0090BR", 78,  "ZC06"),
    @(5272964724, "40AY0090BR", 1,   "ZW04"),
    @(5272964724, "4XE1B81916", 179, "ZC06"),
    @(5272964724, "4X31R64405", 179, "ZC06"),
    @(5272974618, "4XD1P83425", 1,   "ZW08"),
    @(5272970458, "SM10Z35160", 91,  "ZW08"),
    @(5272973084, "4Y41R90027", 112, "ZC06"),
    @(5272973072, "4X21S91185", 20,  "ZW02"),
    @(5272973070, "4X31R64424", 2,   "ZW00"),
    @(5272973068, "4X31R64424", 9,   "ZW08")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}
